## Restored from revision of admin on 03/16/2021 07:57:30 AM.TEST Author: admin. Type: SAVE.
## Change: Rules!C10 18 -> 1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
